$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'lymphedema leggings'
$ws.Range("A2").Value = 'lymphedema tights'
$ws.Range("A3").Value = 'lyra leggings'
$ws.Range("A4").Value = 'lysee leggings'
$ws.Range("A5").Value = 'lysee leggings high waist'
$ws.Range("A6").Value = 'lysee womens leggings'
$ws.Range("A7").Value = 'lysse pants women'
$ws.Range("A8").Value = 'm pro'
$ws.Range("A9").Value = 'm track'
$ws.Range("A10").Value = 'magnolia pants'
$ws.Range("A11").Value = 'make ring tighter'
$ws.Range("A12").Value = 'makenna misuraco'
$ws.Range("A13").Value = 'male compression garment'
$ws.Range("A14").Value = 'male compression tights'
$ws.Range("A15").Value = 'malla de tennis'
$ws.Range("A16").Value = 'malla mujer'
$ws.Range("A17").Value = 'mallas de compresion para mujer'
$ws.Range("A18").Value = 'mallas de futbol'
$ws.Range("A19").Value = 'mallas de licra para mujer'
$ws.Range("A20").Value = 'mallas de mujer'
$ws.Range("A21").Value = 'mallas mujer'
$ws.Range("A22").Value = 'mallas para el gym'
$ws.Range("A23").Value = 'mallas running'
$ws.Range("A24").Value = 'mallas running mujer'
$ws.Range("A25").Value = 'man leggings warm cold weather'
$ws.Range("A26").Value = 'man targets'
$ws.Range("A27").Value = 'marathon 70'
$ws.Range("A28").Value = 'marathon essentials'
$ws.Range("A29").Value = 'marathon frame'
$ws.Range("A30").Value = 'marathon gear'
$ws.Range("A31").Value = 'marathon gear for men'
$ws.Range("A32").Value = 'marathon gear for women'
$ws.Range("A33").Value = 'marathon pants'
$ws.Range("A34").Value = 'marathon runners gear'
$ws.Range("A35").Value = 'marathon running gear'
$ws.Range("A36").Value = 'marathon running tank top women'
$ws.Range("A37").Value = 'marathon running top women'
$ws.Range("A38").Value = 'marathon stick'
$ws.Range("A39").Value = 'marathon tank top'
$ws.Range("A40").Value = 'marathon tank top women'
$ws.Range("A41").Value = 'marathon training for women'
$ws.Range("A42").Value = 'marathon training gear'
$ws.Range("A43").Value = 'marena recovery'
$ws.Range("A44").Value = 'marijuana breath'
$ws.Range("A45").Value = 'marijuana butter'
$ws.Range("A46").Value = 'marijuana crop top for women'
$ws.Range("A47").Value = 'marijuana device'
$ws.Range("A48").Value = 'marijuana drying equipment'
$ws.Range("A49").Value = 'marijuana equipment'
$ws.Range("A50").Value = 'marijuana gift basket'
$ws.Range("A51").Value = 'marijuana jacket'
$ws.Range("A52").Value = 'marijuana joint'
$ws.Range("A53").Value = 'marijuana leggings'
$ws.Range("A54").Value = 'marijuana leggings for women'
$ws.Range("A55").Value = 'marijuana leggins'
$ws.Range("A56").Value = 'marijuana panties'
$ws.Range("A57").Value = 'marijuana penny stocks'
$ws.Range("A58").Value = 'marijuana shorts'
$ws.Range("A59").Value = 'marijuana sweats'
$ws.Range("A60").Value = 'marijuana syndromes'
$ws.Range("A61").Value = 'marijuana training'
$ws.Range("A62").Value = 'marijuana underwear'
$ws.Range("A63").Value = 'marijuana women'
$ws.Range("A64").Value = 'marijuana yoga'
$ws.Range("A65").Value = 'marijuana yoga pants'
$ws.Range("A66").Value = 'mark june'
$ws.Range("A67").Value = 'martial arts ankle support'
$ws.Range("A68").Value = 'martial arts jacket'
$ws.Range("A69").Value = 'martial arts pants'
$ws.Range("A70").Value = 'martial arts ring'
$ws.Range("A71").Value = 'martial arts vest'
$ws.Range("A72").Value = 'master puff ball'
$ws.Range("A73").Value = 'maternity footless tights'
$ws.Range("A74").Value = 'maternity hose'
$ws.Range("A75").Value = 'maternity leggings set'
$ws.Range("A76").Value = 'maternity tights'
$ws.Range("A77").Value = 'maternity workout pants'
$ws.Range("A78").Value = 'maven headbands women'
$ws.Range("A79").Value = 'measurement of joint motion'
$ws.Range("A80").Value = 'media compression'
$ws.Range("A81").Value = 'media de compression mujer'
$ws.Range("A82").Value = 'medias basketball'
$ws.Range("A83").Value = 'medias compresión mujer'
$ws.Range("A84").Value = 'medias compression running'
$ws.Range("A85").Value = 'medias de compresion'
$ws.Range("A86").Value = 'medias de compresión'
$ws.Range("A87").Value = 'medias de compression'
$ws.Range("A88").Value = 'medias de futbol'
$ws.Range("A89").Value = 'medias de malla'
$ws.Range("A90").Value = 'medias de mujer'
$ws.Range("A91").Value = 'medias mallas mujer'
$ws.Range("A92").Value = 'medias sport'
$ws.Range("A93").Value = 'medias women cortas'
$ws.Range("A94").Value = 'medical clothes'
$ws.Range("A95").Value = 'medical clothing'
$ws.Range("A96").Value = 'medical compression capris'
$ws.Range("A97").Value = 'medical compression leggings'
$ws.Range("A98").Value = 'medical compression leggings women'
$ws.Range("A99").Value = 'medical compression pants'
$ws.Range("A100").Value = 'boys running tights'
